# Re-applies the scraped coinranking.com crypto price/volume snapshot
# (and two rank-order swaps: Stacks/PEPE and FirstDigitalUSD/Stellar)
# onto Sheet1. All D (Price) / E (Volume(1h)) cells are plain text in the
# source data (not numbers/percentages), so every write below is forced to
# text with a leading apostrophe and then restored to the "Normal" style so
# Excel does not silently coerce look-alike numbers (e.g. "556.28") into
# real numeric cells or leave a stray quote-prefix style behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextValue "D2" "65.705.62"
Set-TextValue "E2" "  -0.46%  "
Set-TextValue "D3" "3.301.68"
Set-TextValue "E3" "  +0.53%  "
Set-TextValue "E4" "  +0.02%  "
Set-TextValue "D5" "556.28"
Set-TextValue "E5" "  -0.43%  "
Set-TextValue "D6" "184.88"
Set-TextValue "E6" "  -0.53%  "
Set-TextValue "D7" "0.999"
Set-TextValue "E7" "  -0.10%  "
Set-TextValue "D8" "3.293.24"
Set-TextValue "E8" "  +0.48%  "
Set-TextValue "D9" "0.573"
Set-TextValue "E9" "  -3.46%  "
Set-TextValue "D10" "0.174"
Set-TextValue "E10" "  -6.72%  "
Set-TextValue "D11" "0.574"
Set-TextValue "E11" "  -2.18%  "
Set-TextValue "D12" "45.51"
Set-TextValue "E12" "  -4.23%  "
Set-TextValue "D13" "0.0000259"
Set-TextValue "E13" "  -2.61%  "
Set-TextValue "D14" "3.836.32"
Set-TextValue "E14" "  +0.64%  "
Set-TextValue "D15" "8.38"
Set-TextValue "E15" "  -2.96%  "
Set-TextValue "D16" "572.52"
Set-TextValue "E16" "  -10.03%  "
Set-TextValue "D17" "65.692.63"
Set-TextValue "E17" "  -0.44%  "
Set-TextValue "E18" "  +0.30%  "
Set-TextValue "D19" "3.303.44"
Set-TextValue "E19" "  +0.63%  "
Set-TextValue "D20" "17.59"
Set-TextValue "E20" "  -2.34%  "
Set-TextValue "D21" "10.73"
Set-TextValue "E21" "  -5.34%  "
Set-TextValue "D22" "0.886"
Set-TextValue "E22" "  -2.32%  "
Set-TextValue "D23" "17.80"
Set-TextValue "E23" "  -2.91%  "
Set-TextValue "D24" "4.97"
Set-TextValue "E24" "  +1.73%  "
Set-TextValue "D25" "98.17"
Set-TextValue "E25" "  -8.86%  "
Set-TextValue "E26" "  -0.96%  "
Set-TextValue "D27" "2.66"
Set-TextValue "E27" "  -0.36%  "
Set-TextValue "D28" "9.27"
Set-TextValue "E28" "  -3.38%  "
Set-TextValue "D29" "8.42"
Set-TextValue "D30" "30.34"
Set-TextValue "E30" "  +0.00%  "
Set-TextValue "E31" "  +5.69%  "
Set-TextValue "D32" "3.66"
Set-TextValue "E32" "  -8.76%  "
Set-TextValue "D33" "555.29"
Set-TextValue "E33" "  +5.73%  "
Set-TextValue "D34" "10.75"
Set-TextValue "E34" "  -2.63%  "
Set-TextValue "D35" "3.759.45"
Set-TextValue "E35" "  +1.48%  "
Set-TextValue "E36" "  -2.93%  "
Set-TextValue "E37" "  -0.03%  "
Set-TextValue "D38" "55.49"
Set-TextValue "E38" "  -3.67%  "
Set-TextValue "D39" "33.48"
Set-TextValue "E39" "  +1.58%  "
Set-TextValue "E40" "  -2.68%  "
Set-TextValue "B41" "Stacks"
Set-TextValue "C41" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D41" "3.10"
Set-TextValue "E41" "  -7.73%  "
Set-TextValue "B42" "PEPE"
Set-TextValue "C42" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D42" "0.0₃0674"
Set-TextValue "E42" "  -7.53%  "
Set-TextValue "D43" "3.32"
Set-TextValue "E43" "  +0.92%  "
Set-TextValue "D44" "2.55"
Set-TextValue "E44" "  -5.92%  "
Set-TextValue "D45" "0.329"
Set-TextValue "E45" "  -2.89%  "
Set-TextValue "D46" "0.0405"
Set-TextValue "E46" "  -2.21%  "
Set-TextValue "D47" "3.01"
Set-TextValue "E47" "  -7.18%  "
Set-TextValue "B48" "FirstDigitalUSD"
Set-TextValue "C48" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D48" "1.00"
Set-TextValue "E48" "  +0.21%  "
Set-TextValue "B49" "Stellar"
Set-TextValue "C49" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D49" "0.126"
Set-TextValue "E49" "  -2.60%  "
Set-TextValue "E50" "  -3.85%  "
Set-TextValue "D51" "124.92"
Set-TextValue "E51" "  +2.88%  "
